# Slide 9, shape "Rechteck 9" (the "Class SaveEditedRecord" method list):
# rename the logged call from "UtilEmail.sendSecureCallback" to
# "sendNoticationEmailToAdministrator" (notification email sent after the
# record edit), keeping the run's existing formatting intact.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(9)
$sh = $s.Shapes.Item(6)
$tr = $sh.TextFrame.TextRange

$found = $tr.Find("UtilEmail.sendSecureCallback", 0, $false, $false)
if ($found -ne $null) {
    $found.Text = "sendNoticationEmailToAdministrator"
}
